# moved emulators to hyo2_kng + split of SIS listener
#
# - Removes the two old "SIS" icon slides (and their notes pages), keeping
#   the two "Oval" icon slides.
# - Refreshes the cached datetimeFigureOut footer text (2/15/2019 -> 2/17/2019)
#   across the slide master, every slide layout, and the notes master.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText($container, $newText) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

# --- Drop the first two slides (SIS/4 and SIS/5, ids 318 & 575). Their
#     linked notes pages (notesSlide1/2) are removed automatically along
#     with them. The two remaining slides (ids 388 & 574) shift up to
#     become slide 1 and slide 2. ---
$p.Slides.Item(1).Delete()
$p.Slides.Item(1).Delete()

# --- Re-cache the footer date field text on master + every layout. ---
Set-DatePlaceholderText $p.SlideMaster "2/17/2019"

$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Set-DatePlaceholderText $layouts.Item($j) "2/17/2019"
}

# --- Notes master date field uses a different, working code path. ---
$p.NotesMaster.HeadersFooters.DateAndTime.Text = "2/17/2019"

Write-Output ("Slides remaining: " + $p.Slides.Count)
